$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # quality_comparison
$ws2 = $wb.Worksheets.Item(2)   # computational_comparison

# ---------------------------------------------------------------------
# Header borders.
#
# B1:D1 (sheet1) and B1:D1 / E1:G1 (sheet2) are merged header cells.
# Only the anchor cell of a merge (B1 / E1) keeps the bold/centered
# "full box" style; once a border is (re)drawn across the rest of a
# merged range, Excel gives the other physical cells of that range a
# plain, un-fonted style that carries just the visible border edges:
#   - the inner cell(s) get a top+bottom edge
#   - the right-most cell gets top+bottom+right edges
#
# Build the two target looks once (on sheet1 C1 / D1) via explicit
# border edges, then propagate them to every other cell that needs the
# same look with Copy + PasteSpecial(xlPasteFormats) so each cell is
# assigned the very same style index instead of growing new ones.
# ---------------------------------------------------------------------
$xlPasteFormats = -4122

# Build the "top + bottom" look once, from a clean slate.
$innerCell = $ws1.Range("C1")
$innerCell.ClearFormats()
$innerCell.Borders(8).LineStyle = 1    # xlEdgeTop
$innerCell.Borders(9).LineStyle = 1    # xlEdgeBottom

# Derive the "top + bottom + right" look from it (copy, then add the
# one extra edge) instead of building it independently - that keeps
# every border permutation actually visited a "real" final look, so
# no throw-away style/border entries get left behind in styles.xml.
$rightCell = $ws1.Range("D1")
$innerCell.Copy()
$rightCell.PasteSpecial($xlPasteFormats)
$rightCell.Borders(10).LineStyle = 1   # xlEdgeRight

# Propagate both looks to sheet2's matching header cells purely via
# copy/paste of formats so every cell lands on the very same style
# index rather than growing new ones.
$innerCell.Copy()
$ws2.Range("C1").PasteSpecial($xlPasteFormats)
$rightCell.Copy()
$ws2.Range("D1").PasteSpecial($xlPasteFormats)
$innerCell.Copy()
$ws2.Range("F1").PasteSpecial($xlPasteFormats)
$rightCell.Copy()
$ws2.Range("G1").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Anonymize "fedcore" -> "approach" in the header rows of both sheets.
# ---------------------------------------------------------------------
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# ---------------------------------------------------------------------
# Drop the stray empty cell G5 on the computational_comparison sheet.
# ---------------------------------------------------------------------
$ws2.Range("G5").ClearContents()

Write-Host "edit complete"
